# Insert a new weekly record as row 98 on the single data sheet.
# This shifts the existing rows 98-119 down to 99-120 (preserving all of
# their data untouched) and populates the newly inserted row 98 with the
# new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 98; Excel shifts rows 98:119 down to 99:120.
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new record.
$ws.Cells.Item(98, 1).Value = 7
$ws.Cells.Item(98, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(98, 3).Value = "Ñuble"
$ws.Cells.Item(98, 4).Value = 44964
$ws.Cells.Item(98, 5).Value = 16
$ws.Cells.Item(98, 6).Value = 100112031
$ws.Cells.Item(98, 7).Value = "Poroto verde"
$ws.Cells.Item(98, 8).Value = "Sin especificar"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 60
$ws.Cells.Item(98, 11).Value = 26000
$ws.Cells.Item(98, 12).Value = 28000
$ws.Cells.Item(98, 13).Value = 27000
$ws.Cells.Item(98, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(98, 15).Value = "Región del Maule"
$ws.Cells.Item(98, 16).Value = 1080
$ws.Cells.Item(98, 17).Value = 25
$ws.Cells.Item(98, 18).Value = "Hortaliza"
